# Corrijo datos y graficas
# Updates the metrics table (rows 2-6, columns B:M) on the active sheet
# with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(447.0, 76.0, 334.0, 61.0, 0.8546845124282982, 0.8799212598425197, 0.8146341463414634, 0.8671193016488845, 0.8507625272331155, 0.6970151051579176, 0.8466486309913434, 0.9065920875744191)
    3 = @(445.0, 75.0, 335.0, 63.0, 0.8557692307692307, 0.8759842519685039, 0.8170731707317073, 0.8657587548638133, 0.8496732026143791, 0.6950195481771084, 0.846016093385314, 0.9066208949491077)
    4 = @(445.0, 74.0, 336.0, 63.0, 0.8574181117533719, 0.8759842519685039, 0.8195121951219512, 0.866601752677702, 0.8507625272331155, 0.697300966106835, 0.8472778630549538, 0.906808142884579)
    5 = @(445.0, 74.0, 336.0, 63.0, 0.8574181117533719, 0.8759842519685039, 0.8195121951219512, 0.866601752677702, 0.8507625272331155, 0.697300966106835, 0.8472778630549538, 0.906721720760514)
    6 = @(445.0, 74.0, 336.0, 63.0, 0.8574181117533719, 0.8759842519685039, 0.8195121951219512, 0.866601752677702, 0.8507625272331155, 0.697300966106835, 0.8472778630549538, 0.9067169195313992)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
